# Update "想去人数" (F column) figures across sheets to reflect the latest
# generated output (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2907
$ws.Range("F3").Value = 21349
$ws.Range("F5").Value = 3205
$ws.Range("F6").Value = 827
$ws.Range("F13").Value = 133
$ws.Range("F14").Value = 543
$ws.Range("F15").Value = 187
$ws.Range("F16").Value = 303
$ws.Range("F19").Value = 107
$ws.Range("F22").Value = 52

# 演出 (Performances) sheet
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 165

# 本地生活 (Local life) sheet
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6153
$ws.Range("F3").Value = 717
$ws.Range("F4").Value = 714
$ws.Range("F5").Value = 1684
$ws.Range("F6").Value = 68

# 全部类型 (All types) sheet
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6153
$ws.Range("F3").Value = 717
$ws.Range("F4").Value = 714
$ws.Range("F5").Value = 1684
$ws.Range("F6").Value = 2907
$ws.Range("F7").Value = 21349
$ws.Range("F9").Value = 105
$ws.Range("F12").Value = 3206
$ws.Range("F13").Value = 827
$ws.Range("F15").Value = 68
$ws.Range("F25").Value = 133
$ws.Range("F28").Value = 543
$ws.Range("F30").Value = 187
$ws.Range("F32").Value = 303
$ws.Range("F33").Value = 165
$ws.Range("F34").Value = 165
$ws.Range("F38").Value = 107
$ws.Range("F43").Value = 52
